$d = $word.ActiveDocument

# 1) Locate the unique phrase "save it in a drive" so we can target only the
#    standalone word "drive" in the "Unzip it and save it in a drive."
#    sentence (and not the unrelated "drive.google.com" links elsewhere in
#    the document).
$locate = $d.Content
$locate.Find.Execute("save it in a drive", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Narrow the range down to just the word "drive" (the last 5 characters of
# the located phrase) and replace only "rive" with "irectory", leaving the
# initial "d" of "drive" in its own run (mirrors the original authoring edit
# that produced two separate runs: "d" and "irectory").
$driveRange = $d.Range($locate.End - 4, $locate.End)
$driveRange.Find.Execute("rive", $true, $false, $false, $false, $false, $true, 1, $false, "irectory", 2)

# After the replace, $driveRange.End marks the character position right
# after the newly-typed word "directory" (i.e. right before the following
# "."). Remember that position so the bookmark can be re-created there.
$newWordEnd = $driveRange.End

# 2) Remove the existing "_GoBack" bookmark (currently located after the
#    "test examples" text) so it can be re-created at the new location.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 3) Re-insert the "_GoBack" bookmark right after the word "directory" (and
#    before the following period) in the first sentence.
$bookmarkRange = $d.Range($newWordEnd, $newWordEnd)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
